# "Users" workbook - add the Login/Register user table to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 1: Omer Avisror
$ws.Cells.Item(1, 1).Value = 208063511
$ws.Cells.Item(1, 2).Value = "Omer"
$ws.Cells.Item(1, 3).Value = "Avisror"
$ws.Cells.Item(1, 4).Value = 12345
$ws.Cells.Item(1, 5).Value = $true

# Row 2: jimmy james
$ws.Cells.Item(2, 1).Value = 987654321
$ws.Cells.Item(2, 2).Value = "jimmy"
$ws.Cells.Item(2, 3).Value = "james"
$ws.Cells.Item(2, 4).Value = 54321
$ws.Cells.Item(2, 5).Value = $false

# Row 3: Peter Parker
$ws.Cells.Item(3, 1).Value = 123456789
$ws.Cells.Item(3, 2).Value = "Peter"
$ws.Cells.Item(3, 3).Value = "Parker"
$ws.Cells.Item(3, 4).Value = "dibs"
$ws.Cells.Item(3, 5).Value = $true

# Column A is best-fit / custom width to fit the ID numbers.
$ws.Columns.Item(1).ColumnWidth = 9.1

# Leave the selection where the author left it.
$ws.Range("F10").Select() | Out-Null
